$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the rule text for 7 cells so that they describe fricatives as
# "continuant" in addition to "strident" (these cells previously referenced
# the shared strings that are being retired in favor of the new
# "...,continuant,..." variants). The order of assignment below matches the
# order new shared strings are appended to xl/sharedStrings.xml.
$ws.Cells.Item(35,1).Value = "[velar,oral stop]>[coronal,postalveolar,fricative,continuant,posterior,laminal,strident,(backness),(height)]/_[high,front,vowel]"
$ws.Cells.Item(29,1).Value = "[velar,oral stop]>[coronal,alveolar,fricative,continuant,anterior,apical,strident,(backness),(height)]/_[front,vowel]"
$ws.Cells.Item(31,1).Value = "[alveolar,oral stop]>[postalveolar,fricative,continuant,posterior,laminal,strident]/_[front,vowel]"
$ws.Cells.Item(33,1).Value = "[alveolar,oral stop]>[fricative,continuant,strident]/_[front,vowel]"
$ws.Cells.Item(37,1).Value = "[velar,oral stop]>[coronal,alveolar,fricative,continuant,anterior,apical,strident,(backness),(height)]/_[high,front,vowel]"
$ws.Cells.Item(39,1).Value = "[alveolar,oral stop]>[postalveolar,fricative,continuant,posterior,laminal,strident]/_[high,front,vowel]"
$ws.Cells.Item(41,1).Value = "[alveolar,oral stop]>[fricative,continuant,strident]/_[high,front,vowel]"

# Move the visible selection from A27 to A21, reflecting the scrolled-up
# view the author had when saving.
$ws.Range("A21").Select()
